$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows with new allowance ranges, sorted in descending order
$ws.Range("A2").Value = "220.8 - 239.2"
$ws.Range("B2").Value = 222.8
$ws.Range("C2").Value = 227.3
$ws.Range("D2").Value = 223.2
$ws.Range("E2").Value = 227
$ws.Range("F2").Value = 221.6

$ws.Range("A3").Value = "144.0 - 156.0"
$ws.Range("B3").Value = 147.9
$ws.Range("C3").Value = 148
$ws.Range("D3").Value = 149.3
$ws.Range("E3").Value = 148.1
$ws.Range("F3").Value = 150.6

$ws.Range("A4").Value = "67.2 - 72.8"
$ws.Range("B4").Value = 72.8
$ws.Range("C4").Value = 72.8
$ws.Range("D4").Value = 70.90000000000001
$ws.Range("E4").Value = 67.7
$ws.Range("F4").Value = 67.40000000000001
